$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking.com hourly crypto snapshot refresh.
# Most "Price"/"Volume(1h)" cells hold plain text (e.g. "324.85",
# "28.860.49", "  -1.47%  ") even though some of them look numeric.
# Writing a numeric-looking string straight into .Value lets Excel
# auto-coerce it to a real number, which would change the cell type.
# Set-TextValue forces the text number format first (so the literal
# string is kept verbatim) and then restores the cells original
# style, leaving no formatting residue behind.
function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "28.860.49"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "1.877.57"
$ws.Range("E3").Value = "  -2.13%  "

Set-TextValue $ws.Range("D4") "1.004"
$ws.Range("E4").Value = "  +0.04%  "

Set-TextValue $ws.Range("D5") "324.85"
$ws.Range("E5").Value = "  -1.08%  "

$ws.Range("E6").Value = "  -0.01%  "

Set-TextValue $ws.Range("D7") "0.4615"
$ws.Range("E7").Value = "  -1.06%  "

Set-TextValue $ws.Range("D8") "0.3871"
$ws.Range("E8").Value = "  -2.34%  "

Set-TextValue $ws.Range("D9") "0.07839"
$ws.Range("E9").Value = "  -2.53%  "

Set-TextValue $ws.Range("D10") "0.9841"
$ws.Range("E10").Value = "  -3.50%  "

Set-TextValue $ws.Range("D11") "21.73"
$ws.Range("E11").Value = "  -2.71%  "

$ws.Range("D12").Value = "1.879.09"
$ws.Range("E12").Value = "  -1.92%  "

Set-TextValue $ws.Range("D13") "6.988"
$ws.Range("E13").Value = "  -2.39%  "

Set-TextValue $ws.Range("D14") "5.655"
$ws.Range("E14").Value = "  -2.59%  "

Set-TextValue $ws.Range("D15") "0.06977"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("E17").Value = "  +0.10%  "

Set-TextValue $ws.Range("D18") "0.000009952"

Set-TextValue $ws.Range("D19") "16.88"
$ws.Range("E19").Value = "  -2.85%  "

Set-TextValue $ws.Range("D20") "1.003"
$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").Value = "28.859.25"
$ws.Range("E21").Value = "  -1.46%  "

Set-TextValue $ws.Range("D22") "5.251"
$ws.Range("E22").Value = "  -2.56%  "

$ws.Range("E23").Value = "  -1.96%  "

Set-TextValue $ws.Range("D24") "2.104"
$ws.Range("E24").Value = "  +1.83%  "

Set-TextValue $ws.Range("D25") "156.54"

Set-TextValue $ws.Range("D26") "19.31"
$ws.Range("E26").Value = "  -2.44%  "

Set-TextValue $ws.Range("D27") "5.981"
$ws.Range("E27").Value = "  +1.10%  "

Set-TextValue $ws.Range("D28") "117.67"
$ws.Range("E28").Value = "  -2.86%  "

Set-TextValue $ws.Range("D29") "1.907"
$ws.Range("E29").Value = "  -6.35%  "

Set-TextValue $ws.Range("D30") "0.09342"
$ws.Range("E30").Value = "  -0.58%  "

Set-TextValue $ws.Range("D31") "0.9006"
$ws.Range("E31").Value = "  -4.68%  "

Set-TextValue $ws.Range("D32") "5.257"
$ws.Range("E32").Value = "  -2.29%  "

Set-TextValue $ws.Range("D33") "1.317"
$ws.Range("E33").Value = "  -3.65%  "

Set-TextValue $ws.Range("D34") "3.251"
$ws.Range("E34").Value = "  -0.35%  "

$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D35") "1.169"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.05734"
$ws.Range("E36").Value = "  -2.49%  "

Set-TextValue $ws.Range("D37") "0.02071"
$ws.Range("E37").Value = "  -1.92%  "

$ws.Range("E38").Value = "  -0.09%  "

Set-TextValue $ws.Range("D39") "7.638"
$ws.Range("E39").Value = "  -6.50%  "

Set-TextValue $ws.Range("D40") "0.5651"
$ws.Range("E40").Value = "  -3.60%  "

Set-TextValue $ws.Range("D41") "0.1768"
$ws.Range("E41").Value = "  -3.00%  "

Set-TextValue $ws.Range("D42") "9.691"
$ws.Range("E42").Value = "  -4.30%  "

Set-TextValue $ws.Range("D43") "2.231"
$ws.Range("E43").Value = "  -3.89%  "

Set-TextValue $ws.Range("D44") "11.88"
$ws.Range("E44").Value = "  -1.02%  "

Set-TextValue $ws.Range("D45") "0.5334"
$ws.Range("E45").Value = "  -2.76%  "

Set-TextValue $ws.Range("D46") "0.07039"
$ws.Range("E46").Value = "  -2.71%  "

Set-TextValue $ws.Range("D47") "1.837"
$ws.Range("E47").Value = "  -3.05%  "

Set-TextValue $ws.Range("D48") "2.548"
$ws.Range("E48").Value = "  +1.64%  "

Set-TextValue $ws.Range("D49") "112.48"

Set-TextValue $ws.Range("D50") "1.059"
$ws.Range("E50").Value = "  -6.95%  "

Set-TextValue $ws.Range("D51") "70.75"
$ws.Range("E51").Value = "  -1.28%  "
